$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose betting/result data (columns B through AC) were swapped
# between each other in this update, while the leading row index in
# column A stayed with its original row.
$pairs = @(
  @(13,14),
  @(17,18),
  @(19,20),
  @(30,31),
  @(56,57),
  @(74,75),
  @(84,85),
  @(90,91),
  @(106,107)
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($pair in $pairs) {
  $r1 = $pair[0]
  $r2 = $pair[1]

  foreach ($col in $cols) {
    $cell1 = $ws.Range("$col$r1")
    $cell2 = $ws.Range("$col$r2")

    $v1 = $cell1.Value2
    $v2 = $cell2.Value2

    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
  }
}
